$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-03-13 Wednesday" "2024-03-14 Thursday"

Replace-Text "89×18=" "75×21="
Replace-Text "33×89=" "40×23="
Replace-Text "41×65=" "30×11="
Replace-Text "50×80=" "25×29="
Replace-Text "75×49=" "38×30="
Replace-Text "99×80=" "26×43="
Replace-Text "57×16=" "91×92="
Replace-Text "88×20=" "80×66="
Replace-Text "41×25=" "92×51="
Replace-Text "16×61=" "56×95="
Replace-Text "80×87=" "60×80="
Replace-Text "29×45=" "29×56="
Replace-Text "83×48=" "85×37="
Replace-Text "20×85=" "22×87="
Replace-Text "76×92=" "91×62="
Replace-Text "22×45=" "79×26="
Replace-Text "62×69=" "22×55="
Replace-Text "18×25=" "95×80="
Replace-Text "73×46=" "53×48="
Replace-Text "35×65=" "52×94="
Replace-Text "18×52=" "35×97="
Replace-Text "70×24=" "24×69="
Replace-Text "32×94=" "13×75="
Replace-Text "42×69=" "13×28="
Replace-Text "63×58=" "39×69="
